$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 5 - España
$ws.Range("B5").Value = 229422
$ws.Range("C5").Value = 2793
$ws.Range("D5").Value = 120832

# Row 8 - Alemania
$ws.Range("B8").Value = 157946
$ws.Range("C8").Value = 176
$ws.Range("E8").Value = 37462
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 5984

# Row 17 - Paises Bajos
$ws.Range("B17").Value = 38245
$ws.Range("C17").Value = 400
$ws.Range("E17").Value = 33477
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = 4518

# Row 18 - Suiza
$ws.Range("B18").Value = 29164
$ws.Range("C18").Value = 103
$ws.Range("E18").Value = 5724

# Row 21 - Portugal
$ws.Range("B21").Value = 24027
$ws.Range("C21").Value = 163
$ws.Range("D21").Value = 1357
$ws.Range("E21").Value = 21742
$ws.Range("F21").Value = 176
$ws.Range("G21").Value = 25
$ws.Range("H21").Value = 928

# Row 24 - Suecia
$ws.Range("B24").Value = 18926
$ws.Range("C24").Value = 286
$ws.Range("E24").Value = 15647
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = 2274

# Row 67 - Croacia
$ws.Range("B67").Value = 2039
$ws.Range("C67").Value = 9
$ws.Range("D67").Value = 1166
$ws.Range("E67").Value = 814
$ws.Range("F67").Value = 21
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 59

# Row 68 - Uzbekistan
$ws.Range("D68").Value = 880
$ws.Range("E68").Value = 999
